# Add base Frame work / small test
# Reproduces the changes made to Scene.xlsx:
#  - two new rows of data (Menu / MainGame) appended to the table
#  - the whole A1:E6 table left-aligned
#  - selection / dimension updated to match the new data extent

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Add the two new data rows
# ---------------------------------------------------------------------
$ws.Range("B5").Value = 1
$ws.Range("D5").Value = "Menu"
$ws.Range("E5").Value = 0

$ws.Range("B6").Value = 2
$ws.Range("D6").Value = "MainGame"
$ws.Range("E6").Value = 0

# ---------------------------------------------------------------------
# 2. Left align every populated cell in the A1:E6 block (matches the
#    "horizontal=left" alignment that was added to every cell style).
#    Cells that used the "Scene ID / Resource Name" font are re-stamped
#    with their font explicitly first so the font is not lost when the
#    alignment is (re)applied.
# ---------------------------------------------------------------------
$fontCells = @("B4", "D4", "D5", "D6")
foreach ($addr in $fontCells) {
    $ws.Range($addr).Font.Name = "ＭＳ Ｐゴシック"
    $ws.Range($addr).Font.Size = 11
}

$allCells = @(
    "A1","B1","C1","D1","E1",
    "A2","B2","C2","D2","E2",
    "A3","B3","C3","D3","E3",
    "A4","B4","C4","D4","E4",
    "A5","B5","C5","D5","E5",
    "A6","B6","C6","D6","E6"
)
foreach ($addr in $allCells) {
    $ws.Range($addr).HorizontalAlignment = -4131
}

# ---------------------------------------------------------------------
# 3. Update the selection to the full used range
# ---------------------------------------------------------------------
$ws.Range("A1:E6").Select()

Write-Output "Scene.xlsx updated"
